$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update shared-string / rich-text header cells ---
$a8 = $ws.Range("A8")
$a8.Value = $a8.Value().Replace("32", "33")

$c9 = $ws.Range("C9")
$c9.Value = $c9.Value().Replace("8/7/2023", "8/14/2023").Replace("8/13/2023", "8/20/2023")

# --- Update numeric data cells (rows 14-30) ---
$ws.Range("D14").Value = 1
$ws.Range("J14").Value = 34
$ws.Range("K14").Value = -67.647058823529
$ws.Range("L14").Value = -64.516129032258
$ws.Range("M14").Value = -75
$ws.Range("N14").Value = -87.640449438202
$ws.Range("C15").Value = 2
$ws.Range("D15").Value = 4
$ws.Range("E15").Value = -50
$ws.Range("G15").Value = 13
$ws.Range("H15").Value = -23.076923076923
$ws.Range("I15").Value = 90
$ws.Range("J15").Value = 123
$ws.Range("K15").Value = -26.829268292682
$ws.Range("L15").Value = -21.052631578947
$ws.Range("M15").Value = 1.123595505617
$ws.Range("N15").Value = -60.176991150442
$ws.Range("C16").Value = 33
$ws.Range("D16").Value = 36
$ws.Range("E16").Value = -8.333333333333
$ws.Range("F16").Value = 117
$ws.Range("G16").Value = 139
$ws.Range("H16").Value = -15.827338129496
$ws.Range("I16").Value = 844
$ws.Range("J16").Value = 945
$ws.Range("K16").Value = -10.687830687830
$ws.Range("L16").Value = 22.318840579710
$ws.Range("M16").Value = -36.445783132530
$ws.Range("N16").Value = -82.449573715949
$ws.Range("C17").Value = 68
$ws.Range("D17").Value = 49
$ws.Range("E17").Value = 38.775510204081
$ws.Range("F17").Value = 236
$ws.Range("G17").Value = 206
$ws.Range("H17").Value = 14.563106796116
$ws.Range("I17").Value = 1780
$ws.Range("J17").Value = 1751
$ws.Range("K17").Value = 1.656196459166
$ws.Range("L17").Value = 14.543114543114
$ws.Range("M17").Value = 63.602941176470
$ws.Range("N17").Value = -24.287537218205
$ws.Range("C18").Value = 20
$ws.Range("D18").Value = 23
$ws.Range("E18").Value = -13.043478260869
$ws.Range("F18").Value = 94
$ws.Range("G18").Value = 99
$ws.Range("H18").Value = -5.050505050505
$ws.Range("I18").Value = 725
$ws.Range("J18").Value = 716
$ws.Range("K18").Value = 1.256983240223
$ws.Range("L18").Value = 26.086956521739
$ws.Range("M18").Value = -44.486983154670
$ws.Range("N18").Value = -87.512917671374
$ws.Range("C19").Value = 84
$ws.Range("D19").Value = 74
$ws.Range("E19").Value = 13.513513513513
$ws.Range("F19").Value = 319
$ws.Range("G19").Value = 304
$ws.Range("H19").Value = 4.934210526315
$ws.Range("I19").Value = 2260
$ws.Range("J19").Value = 2384
$ws.Range("K19").Value = -5.201342281879
$ws.Range("L19").Value = 42.138364779874
$ws.Range("M19").Value = 29.810453762205
$ws.Range("N19").Value = -56.771231828615
$ws.Range("C20").Value = 35
$ws.Range("D20").Value = 19
$ws.Range("E20").Value = 84.210526315789
$ws.Range("F20").Value = 164
$ws.Range("G20").Value = 96
$ws.Range("H20").Value = 70.833333333333
$ws.Range("I20").Value = 1125
$ws.Range("J20").Value = 978
$ws.Range("K20").Value = 15.030674846625
$ws.Range("L20").Value = 68.918918918918
$ws.Range("M20").Value = 4.748603351955
$ws.Range("N20").Value = -90.627343164209
$ws.Range("C21").Value = 242
$ws.Range("D21").Value = 206
$ws.Range("E21").Value = 17.475728155339
$ws.Range("F21").Value = 940
$ws.Range("G21").Value = 865
$ws.Range("H21").Value = 8.670520231213
$ws.Range("I21").Value = 6835
$ws.Range("J21").Value = 6931
$ws.Range("K21").Value = -1.385081517818
$ws.Range("L21").Value = 30.938697318007
$ws.Range("M21").Value = 2.473763118440
$ws.Range("N21").Value = -77.598977451494
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = 33.333333333333
$ws.Range("F22").Value = 10
$ws.Range("G22").Value = 15
$ws.Range("H22").Value = -33.333333333333
$ws.Range("I22").Value = 75
$ws.Range("J22").Value = 73
$ws.Range("K22").Value = 2.739726027397
$ws.Range("L22").Value = 17.1875
$ws.Range("M22").Value = -5.063291139240
$ws.Range("C23").Value = 8
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 22
$ws.Range("G23").Value = 21
$ws.Range("H23").Value = 4.761904761904
$ws.Range("I23").Value = 160
$ws.Range("J23").Value = 145
$ws.Range("K23").Value = 10.344827586206
$ws.Range("L23").Value = 13.475177304964
$ws.Range("M23").Value = 68.421052631578
$ws.Range("C24").Value = 173
$ws.Range("D24").Value = 231
$ws.Range("E24").Value = -25.108225108225
$ws.Range("F24").Value = 685
$ws.Range("G24").Value = 866
$ws.Range("H24").Value = -20.900692840646
$ws.Range("I24").Value = 5896
$ws.Range("J24").Value = 6293
$ws.Range("K24").Value = -6.308596853646
$ws.Range("L24").Value = 40.917782026768
$ws.Range("M24").Value = 50.94726062468
$ws.Range("C25").Value = 76
$ws.Range("D25").Value = 71
$ws.Range("E25").Value = 7.042253521126
$ws.Range("F25").Value = 370
$ws.Range("G25").Value = 337
$ws.Range("H25").Value = 9.792284866468
$ws.Range("I25").Value = 2877
$ws.Range("J25").Value = 2584
$ws.Range("K25").Value = 11.339009287925
$ws.Range("L25").Value = 32.702952029520
$ws.Range("M25").Value = -2.869682646860
$ws.Range("C26").Value = 2
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -66.666666666666
$ws.Range("F26").Value = 14
$ws.Range("G26").Value = 19
$ws.Range("H26").Value = -26.315789473684
$ws.Range("I26").Value = 160
$ws.Range("J26").Value = 195
$ws.Range("K26").Value = -17.948717948717
$ws.Range("L26").Value = -9.090909090909
$ws.Range("C27").Value = 6
$ws.Range("D27").Value = 8
$ws.Range("E27").Value = -25
$ws.Range("F27").Value = 32
$ws.Range("G27").Value = 36
$ws.Range("H27").Value = -11.111111111111
$ws.Range("I27").Value = 271
$ws.Range("J27").Value = 278
$ws.Range("K27").Value = -2.517985611510
$ws.Range("L27").Value = 7.539682539682
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = -75
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 21
$ws.Range("H28").Value = -85.714285714285
$ws.Range("I28").Value = 59
$ws.Range("J28").Value = 126
$ws.Range("K28").Value = -53.174603174603
$ws.Range("L28").Value = -60.135135135135
$ws.Range("M28").Value = -54.615384615384
$ws.Range("N28").Value = -81.028938906752
$ws.Range("C29").Value = 1
$ws.Range("E29").Value = -75
$ws.Range("F29").Value = 3
$ws.Range("G29").Value = 16
$ws.Range("H29").Value = -81.25
$ws.Range("I29").Value = 43
$ws.Range("J29").Value = 96
$ws.Range("K29").Value = -55.208333333333
$ws.Range("L29").Value = -64.462809917355
$ws.Range("M29").Value = -59.047619047619
$ws.Range("N29").Value = -84.965034965035
$ws.Range("D30").Value = 2
$ws.Range("G30").Value = 5
$ws.Range("J30").Value = 25
$ws.Range("K30").Value = 36

# --- Fix style for C28/C29 (was text "0", now numeric) ---
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("C29").NumberFormat = "#,##0"
